# "Ready Labs and products data"
#
# 1. Rename the "optimals" sheet to "ideals".
# 2. Make the "ideals" sheet the active/selected tab (it was previously
#    on "ideales"), which also resets its view back to the top.
# 3. Rename the builtin "Excel Built-in Comma [0]" cell style to
#    "Excel Built-in Explanatory Text".

$wb = $excel.ActiveWorkbook

# 1. Rename "optimals" -> "ideals"
$ws = $wb.Worksheets.Item("optimals")
$ws.Name = "ideals"

# 2. Activate the renamed sheet so it becomes the selected/visible tab
$ws.Activate()

# 3. Rename the builtin cell style
$style = $wb.Styles.Item("Excel Built-in Comma [0]")
$style.Name = "Excel Built-in Explanatory Text"
